$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Set Values Here")

# Row 8 corresponds to "carbon tax revenue" weighting row.
# Update the government-revenue-accounting weights for the carbon tax:
#   Deficit Spending (C8): 0 -> 5
#   Payroll Taxes   (E8): 0 -> 5
$ws.Range("C8").Value = 5
$ws.Range("E8").Value = 5

$wb.Save()
